$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186012625694275
$ws.Range("B1").Value = 2.269022703170776
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.333497047424316
$ws.Range("E1").Value = 1.218836188316345
